# Applies the cell-value updates described in the commit diff to the
# "Anima_Profits" workbook (profit/cost calculation columns H..N across
# the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets).
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
#          J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ,
#          M=LeveProfitNQ, N=LeveProfitHQ
#
# Setting a cell's .Value to $null removes it entirely (matching rows in
# the diff where a M/N cell disappears); conversely, assigning a numeric
# value to a previously-empty cell creates it.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 44
$ws.Cells.Item(44, 8).Value = 5000
$ws.Cells.Item(44, 10).Value = 5000
$ws.Cells.Item(44, 12).Value = 5000
$ws.Cells.Item(44, 14).Value = -5924
# Row 51
$ws.Cells.Item(51, 8).Value = 1209.1428
$ws.Cells.Item(51, 9).Value = 988
$ws.Cells.Item(51, 10).Value = 1375
$ws.Cells.Item(51, 11).Value = 988
$ws.Cells.Item(51, 12).Value = 1375
$ws.Cells.Item(51, 13).Value = -504
$ws.Cells.Item(51, 14).Value = -2343
# Row 64
$ws.Cells.Item(64, 8).Value = 3016.647
$ws.Cells.Item(64, 9).Value = 2992.2
$ws.Cells.Item(64, 11).Value = 2992.2
$ws.Cells.Item(64, 13).Value = -2744.2
# Row 67
$ws.Cells.Item(67, 8).Value = 3016.647
$ws.Cells.Item(67, 9).Value = 2992.2
$ws.Cells.Item(67, 11).Value = 2992.2
$ws.Cells.Item(67, 13).Value = -2134.2
# Row 70
$ws.Cells.Item(70, 8).Value = 1480.8182
$ws.Cells.Item(70, 9).Value = 1333.3334
$ws.Cells.Item(70, 10).Value = 1536.125
$ws.Cells.Item(70, 11).Value = 4000.0002
$ws.Cells.Item(70, 12).Value = 4608.375
$ws.Cells.Item(70, 13).Value = -3730.0002
$ws.Cells.Item(70, 14).Value = -5148.375
# Row 73
$ws.Cells.Item(73, 8).Value = 1480.8182
$ws.Cells.Item(73, 9).Value = 1333.3334
$ws.Cells.Item(73, 10).Value = 1536.125
$ws.Cells.Item(73, 11).Value = 4000.0002
$ws.Cells.Item(73, 12).Value = 4608.375
$ws.Cells.Item(73, 13).Value = -3064.0002
$ws.Cells.Item(73, 14).Value = -6480.375
# Row 137
$ws.Cells.Item(137, 8).Value = 1135.7778
$ws.Cells.Item(137, 9).Value = 526.2727
$ws.Cells.Item(137, 11).Value = 1578.8181
$ws.Cells.Item(137, 13).Value = 971.1819
# Row 138
$ws.Cells.Item(138, 8).Value = 3122.6191
$ws.Cells.Item(138, 9).Value = 3251.6365
$ws.Cells.Item(138, 10).Value = 3076.8386
$ws.Cells.Item(138, 11).Value = 9754.9095
$ws.Cells.Item(138, 12).Value = 9230.515800000001
$ws.Cells.Item(138, 13).Value = -4614.9095
$ws.Cells.Item(138, 14).Value = -19510.5158

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 17650.334
$ws.Cells.Item(2, 9).Value = 33993.332
$ws.Cells.Item(2, 10).Value = 1307.3334
$ws.Cells.Item(2, 11).Value = 33993.332
$ws.Cells.Item(2, 12).Value = 1307.3334
$ws.Cells.Item(2, 13).Value = -33880.332
$ws.Cells.Item(2, 14).Value = -1533.3334
# Row 105
$ws.Cells.Item(105, 8).Value = 73123.336
$ws.Cells.Item(105, 10).Value = 73123.336
$ws.Cells.Item(105, 12).Value = 73123.336
$ws.Cells.Item(105, 14).Value = -80111.336
# Row 106
$ws.Cells.Item(106, 8).Value = 44444
$ws.Cells.Item(106, 10).Value = 44444
$ws.Cells.Item(106, 12).Value = 44444
$ws.Cells.Item(106, 14).Value = -46968
# Row 116
$ws.Cells.Item(116, 8).Value = 17650.334
$ws.Cells.Item(116, 9).Value = 33993.332
$ws.Cells.Item(116, 10).Value = 1307.3334
$ws.Cells.Item(116, 11).Value = 33993.332
$ws.Cells.Item(116, 12).Value = 1307.3334
$ws.Cells.Item(116, 13).Value = -31699.332
$ws.Cells.Item(116, 14).Value = -5895.3334

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 17650.334
$ws.Cells.Item(3, 9).Value = 33993.332
$ws.Cells.Item(3, 10).Value = 1307.3334
$ws.Cells.Item(3, 11).Value = 33993.332
$ws.Cells.Item(3, 12).Value = 1307.3334
$ws.Cells.Item(3, 13).Value = -33879.332
$ws.Cells.Item(3, 14).Value = -1535.3334
# Row 134
$ws.Cells.Item(134, 8).Value = 2944
$ws.Cells.Item(134, 9).Value = 3053.739
$ws.Cells.Item(134, 10).Value = 2733.6667
$ws.Cells.Item(134, 11).Value = 9161.217000000001
$ws.Cells.Item(134, 12).Value = 8201.000100000001
$ws.Cells.Item(134, 13).Value = -6626.217000000001
$ws.Cells.Item(134, 14).Value = -13271.0001

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Cells.Item(99, 8).Value = 1504.3334
$ws.Cells.Item(99, 9).Value = 1499
$ws.Cells.Item(99, 10).Value = 1507
$ws.Cells.Item(99, 11).Value = 1499
$ws.Cells.Item(99, 12).Value = 1507
$ws.Cells.Item(99, 13).Value = -1
$ws.Cells.Item(99, 14).Value = -4503
# Row 112
$ws.Cells.Item(112, 8).Value = 40702
$ws.Cells.Item(112, 10).Value = 40702
$ws.Cells.Item(112, 12).Value = 40702
$ws.Cells.Item(112, 14).Value = -43656
# Row 122
$ws.Cells.Item(122, 8).Value = 1448.8334
$ws.Cells.Item(122, 9).Value = 1448.8334
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4346.5002
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1896.5002
$ws.Cells.Item(122, 14).Value = $null
# Row 126
$ws.Cells.Item(126, 8).Value = 1504.3334
$ws.Cells.Item(126, 9).Value = 1499
$ws.Cells.Item(126, 10).Value = 1507
$ws.Cells.Item(126, 11).Value = 4497
$ws.Cells.Item(126, 12).Value = 4521
$ws.Cells.Item(126, 13).Value = -2027
$ws.Cells.Item(126, 14).Value = -9461
# Row 135
$ws.Cells.Item(135, 8).Value = 50210
$ws.Cells.Item(135, 10).Value = 50210
$ws.Cells.Item(135, 12).Value = 50210
$ws.Cells.Item(135, 14).Value = -60350

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Cells.Item(12, 8).Value = 172
$ws.Cells.Item(12, 9).Value = 183.33333
$ws.Cells.Item(12, 10).Value = 163.5
$ws.Cells.Item(12, 11).Value = 549.99999
$ws.Cells.Item(12, 12).Value = 490.5
$ws.Cells.Item(12, 13).Value = -376.99999
$ws.Cells.Item(12, 14).Value = -836.5
# Row 107
$ws.Cells.Item(107, 8).Value = 1421.965
$ws.Cells.Item(107, 10).Value = 2323.4375
$ws.Cells.Item(107, 12).Value = 6970.3125
$ws.Cells.Item(107, 14).Value = -10810.3125
# Row 137
$ws.Cells.Item(137, 8).Value = 8255.053
$ws.Cells.Item(137, 10).Value = 3450
$ws.Cells.Item(137, 12).Value = 10350
$ws.Cells.Item(137, 14).Value = -20550

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 1347.8511
$ws.Cells.Item(68, 9).Value = 1265.1666
$ws.Cells.Item(68, 11).Value = 1265.1666
$ws.Cells.Item(68, 13).Value = -516.1666
# Row 71
$ws.Cells.Item(71, 8).Value = 1347.8511
$ws.Cells.Item(71, 9).Value = 1265.1666
$ws.Cells.Item(71, 11).Value = 6325.833000000001
$ws.Cells.Item(71, 13).Value = -2581.833000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Cells.Item(21, 8).Value = 70017
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 70017
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 70017
$ws.Cells.Item(21, 13).Value = $null
$ws.Cells.Item(21, 14).Value = -70487
# Row 33
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 14).Value = $null
# Row 34
$ws.Cells.Item(34, 8).Value = 40000
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).Value = $null
# Row 35
$ws.Cells.Item(35, 8).Value = 70017
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 70017
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 70017
$ws.Cells.Item(35, 13).Value = $null
$ws.Cells.Item(35, 14).Value = -70597
# Row 36
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 14).Value = $null
# Row 37
$ws.Cells.Item(37, 8).Value = 62521.75
$ws.Cells.Item(37, 10).Value = 62521.75
$ws.Cells.Item(37, 12).Value = 62521.75
$ws.Cells.Item(37, 14).Value = -62927.75
# Row 39
$ws.Cells.Item(39, 8).Value = 40000
$ws.Cells.Item(39, 10).Value = 40000
$ws.Cells.Item(39, 12).Value = 40000
$ws.Cells.Item(39, 14).Value = -40826
# Row 40
$ws.Cells.Item(40, 8).Value = 40000
$ws.Cells.Item(40, 10).Value = 40000
$ws.Cells.Item(40, 12).Value = 40000
$ws.Cells.Item(40, 14).Value = -40298
# Row 41
$ws.Cells.Item(41, 8).Value = 7500
$ws.Cells.Item(41, 10).Value = 7500
$ws.Cells.Item(41, 12).Value = 7500
$ws.Cells.Item(41, 14).Value = -8280
# Row 45
$ws.Cells.Item(45, 8).Value = 9392
$ws.Cells.Item(45, 9).Value = 8189.3335
$ws.Cells.Item(45, 10).Value = 13000
$ws.Cells.Item(45, 11).Value = 8189.3335
$ws.Cells.Item(45, 12).Value = 13000
$ws.Cells.Item(45, 13).Value = -7698.3335
$ws.Cells.Item(45, 14).Value = -13982
# Row 52
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).Value = $null
# Row 59
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).Value = $null
# Row 62
$ws.Cells.Item(62, 8).Value = 40799
$ws.Cells.Item(62, 9).Value = 3823
$ws.Cells.Item(62, 10).Value = 77775
$ws.Cells.Item(62, 11).Value = 3823
$ws.Cells.Item(62, 12).Value = 77775
$ws.Cells.Item(62, 13).Value = -3199
$ws.Cells.Item(62, 14).Value = -79023
# Row 63
$ws.Cells.Item(63, 8).Value = 89624.5
$ws.Cells.Item(63, 10).Value = 89624.5
$ws.Cells.Item(63, 12).Value = 89624.5
$ws.Cells.Item(63, 14).Value = -90872.5
# Row 65
$ws.Cells.Item(65, 8).Value = 40799
$ws.Cells.Item(65, 9).Value = 3823
$ws.Cells.Item(65, 10).Value = 77775
$ws.Cells.Item(65, 11).Value = 19115
$ws.Cells.Item(65, 12).Value = 388875
$ws.Cells.Item(65, 13).Value = -15995
$ws.Cells.Item(65, 14).Value = -395115
# Row 66
$ws.Cells.Item(66, 8).Value = 89624.5
$ws.Cells.Item(66, 10).Value = 89624.5
$ws.Cells.Item(66, 12).Value = 268873.5
$ws.Cells.Item(66, 14).Value = -275113.5
# Row 74
$ws.Cells.Item(74, 8).Value = 4617.3335
$ws.Cells.Item(74, 10).Value = 4617.3335
$ws.Cells.Item(74, 12).Value = 4617.3335
$ws.Cells.Item(74, 14).Value = -6489.3335
# Row 77
$ws.Cells.Item(77, 8).Value = 4617.3335
$ws.Cells.Item(77, 10).Value = 4617.3335
$ws.Cells.Item(77, 12).Value = 13852.0005
$ws.Cells.Item(77, 14).Value = -23212.0005
# Row 120
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 14).Value = $null
# Row 121
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 14).Value = $null
# Row 123
$ws.Cells.Item(123, 8).Value = 28000
$ws.Cells.Item(123, 10).Value = 28000
$ws.Cells.Item(123, 12).Value = 28000
$ws.Cells.Item(123, 14).Value = -37800
# Row 124
$ws.Cells.Item(124, 8).Value = 25607.125
$ws.Cells.Item(124, 10).Value = 25607.125
$ws.Cells.Item(124, 12).Value = 25607.125
$ws.Cells.Item(124, 14).Value = -35427.125
# Row 125
$ws.Cells.Item(125, 8).Value = 10715
$ws.Cells.Item(125, 10).Value = 10715
$ws.Cells.Item(125, 12).Value = 10715
$ws.Cells.Item(125, 14).Value = -20555
# Row 127
$ws.Cells.Item(127, 8).Value = 64178.625
$ws.Cells.Item(127, 10).Value = 64178.625
$ws.Cells.Item(127, 12).Value = 64178.625
$ws.Cells.Item(127, 14).Value = -74098.625
# Row 130
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).Value = $null
# Row 131
$ws.Cells.Item(131, 8).Value = 55357.5
$ws.Cells.Item(131, 9).Value = 10000
$ws.Cells.Item(131, 10).Value = 100715
$ws.Cells.Item(131, 11).Value = 10000
$ws.Cells.Item(131, 12).Value = 100715
$ws.Cells.Item(131, 13).Value = -4960
$ws.Cells.Item(131, 14).Value = -110795
# Row 132
$ws.Cells.Item(132, 8).Value = 6207177.5
$ws.Cells.Item(132, 9).Value = 1369.8975
$ws.Cells.Item(132, 11).Value = 4109.6925
$ws.Cells.Item(132, 13).Value = -1579.6925
# Row 135
$ws.Cells.Item(135, 8).Value = 62266.54
$ws.Cells.Item(135, 10).Value = 62266.54
$ws.Cells.Item(135, 12).Value = 62266.54
$ws.Cells.Item(135, 14).Value = -72406.54000000001